$d = $word.ActiveDocument

$replacements = @(
    @("809×9=", "154×2="),
    @("778×4=", "746×5="),
    @("949×4=", "782×8="),
    @("111×8=", "576×5="),
    @("109×9=", "912×8="),
    @("877×5=", "448×5="),
    @("114×7=", "640×7="),
    @("315×9=", "922×5="),
    @("982×9=", "607×4="),
    @("552×6=", "967×8="),
    @("786×2=", "961×5="),
    @("396×8=", "843×3="),
    @("442×9=", "278×3="),
    @("456×4=", "351×4="),
    @("365×4=", "627×3="),
    @("855×3=", "186×6="),
    @("992×2=", "748×9="),
    @("619×7=", "581×9="),
    @("955×2=", "401×2="),
    @("819×6=", "746×2="),
    @("912×7=", "178×6="),
    @("842×3=", "243×4="),
    @("925×4=", "597×7="),
    @("313×6=", "107×5="),
    @("908×8=", "263×7=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
